# Update "想去人数" (want-to-go count) figures for two events on both the
# "展览" and "全部类型" worksheets, matching the refreshed data snapshot.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 7430
    $ws.Range("F3").Value = 7417
    $ws.Range("F17").Value = 3
}
